$wb = $excel.ActiveWorkbook

# ---- Table 1 (Sheet1): reorder rows, drop "complete: n = ..." lines ----
$ws1 = $wb.Worksheets.Item(1)

$ws1.Cells.Item(2, 1).Value = "Sex"
$ws1.Cells.Item(2, 2).Value = "female: 67% (n = 320)`nmale: 33% (n = 159)"
$ws1.Cells.Item(2, 3).Value = "female: 70% (n = 300)`nmale: 30% (n = 127)"
$ws1.Cells.Item(2, 4).Value = "ns (p = 0.46)"
$ws1.Cells.Item(2, 5).Value = "V = 0.037"

$ws1.Cells.Item(3, 1).Value = "Education"
$ws1.Cells.Item(3, 2).Value = "non-tertiary: 63% (n = 302)`ntertiary: 37% (n = 176)"
$ws1.Cells.Item(3, 3).Value = "non-tertiary: 59% (n = 250)`ntertiary: 41% (n = 177)"
$ws1.Cells.Item(3, 4).Value = "ns (p = 0.35)"
$ws1.Cells.Item(3, 5).Value = "V = 0.047"

$ws1.Cells.Item(4, 1).Value = "Age, years"
$ws1.Cells.Item(4, 2).Value = "43 [IQR: 32 - 53]`nrange: 18 - 80"
$ws1.Cells.Item(4, 3).Value = "45 [IQR: 34 - 54]`nrange: 18 - 95"
$ws1.Cells.Item(4, 4).Value = "ns (p = 0.31)"
$ws1.Cells.Item(4, 5).Value = "r = 0.048"

$ws1.Cells.Item(5, 1).Value = "BMI before COVID-19"
$ws1.Cells.Item(5, 2).Value = "normal: 54% (n = 257)`noverweight: 28% (n = 135)`nobesity: 18% (n = 84)"
$ws1.Cells.Item(5, 3).Value = "normal: 66% (n = 278)`noverweight: 25% (n = 104)`nobesity: 8.8% (n = 37)"
$ws1.Cells.Item(5, 4).Value = "p = 0.0011"
$ws1.Cells.Item(5, 5).Value = "V = 0.15"

$ws1.Cells.Item(6, 1).Value = "Employment status"
$ws1.Cells.Item(6, 2).Value = "employed: 83% (n = 398)`nunemployed: 8.4% (n = 40)`nleave: 1.7% (n = 8)`nretired: 6.9% (n = 33)"
$ws1.Cells.Item(6, 3).Value = "employed: 81% (n = 348)`nunemployed: 9.4% (n = 40)`nleave: 1.9% (n = 8)`nretired: 7.3% (n = 31)"
$ws1.Cells.Item(6, 4).Value = "ns (p = 1)"
$ws1.Cells.Item(6, 5).Value = "V = 0.022"

$ws1.Cells.Item(7, 1).Value = "Autoimmunity"
$ws1.Cells.Item(7, 2).Value = "6.7% (n = 32)"
$ws1.Cells.Item(7, 3).Value = "6.3% (n = 27)"
$ws1.Cells.Item(7, 4).Value = "ns (p = 1)"
$ws1.Cells.Item(7, 5).Value = "V = 0.0072"

$ws1.Cells.Item(8, 1).Value = "Hypertension"
$ws1.Cells.Item(8, 2).Value = "11% (n = 51)"
$ws1.Cells.Item(8, 3).Value = "8.4% (n = 36)"
$ws1.Cells.Item(8, 4).Value = "ns (p = 0.46)"
$ws1.Cells.Item(8, 5).Value = "V = 0.038"

$ws1.Cells.Item(9, 1).Value = "Pre-CoV depression/anxiety"
$ws1.Cells.Item(9, 2).Value = "5.4% (n = 26)"
$ws1.Cells.Item(9, 3).Value = "5.2% (n = 22)"
$ws1.Cells.Item(9, 4).Value = "ns (p = 1)"
$ws1.Cells.Item(9, 5).Value = "V = 0.0061"

$ws1.Cells.Item(10, 1).Value = "Diabetes"
$ws1.Cells.Item(10, 2).Value = "1.5% (n = 7)"
$ws1.Cells.Item(10, 3).Value = "0.23% (n = 1)"
$ws1.Cells.Item(10, 4).Value = "ns (p = 0.26)"
$ws1.Cells.Item(10, 5).Value = "V = 0.065"

$ws1.Cells.Item(11, 1).Value = "Freq. resp. infections"
$ws1.Cells.Item(11, 2).Value = "6.7% (n = 32)"
$ws1.Cells.Item(11, 3).Value = "3.3% (n = 14)"
$ws1.Cells.Item(11, 4).Value = "ns (p = 0.1)"
$ws1.Cells.Item(11, 5).Value = "V = 0.077"

$ws1.Cells.Item(12, 1).Value = "Cardiovascular disease"
$ws1.Cells.Item(12, 2).Value = "2.1% (n = 10)"
$ws1.Cells.Item(12, 3).Value = "3% (n = 13)"
$ws1.Cells.Item(12, 4).Value = "ns (p = 0.62)"
$ws1.Cells.Item(12, 5).Value = "V = 0.03"

$ws1.Cells.Item(13, 1).Value = "Hay fever/allergy"
$ws1.Cells.Item(13, 2).Value = "18% (n = 88)"
$ws1.Cells.Item(13, 3).Value = "12% (n = 51)"
$ws1.Cells.Item(13, 4).Value = "p = 0.045"
$ws1.Cells.Item(13, 5).Value = "V = 0.089"

$ws1.Cells.Item(14, 1).Value = "Malignancy"
$ws1.Cells.Item(14, 2).Value = "2.1% (n = 10)"
$ws1.Cells.Item(14, 3).Value = "4% (n = 17)"
$ws1.Cells.Item(14, 4).Value = "ns (p = 0.31)"
$ws1.Cells.Item(14, 5).Value = "V = 0.056"

$ws1.Cells.Item(15, 1).Value = "Gastrointestinal disease"
$ws1.Cells.Item(15, 2).Value = "1.7% (n = 8)"
$ws1.Cells.Item(15, 3).Value = "0.7% (n = 3)"
$ws1.Cells.Item(15, 4).Value = "ns (p = 0.46)"
$ws1.Cells.Item(15, 5).Value = "V = 0.044"

$ws1.Cells.Item(16, 1).Value = "Pulmonary disease"
$ws1.Cells.Item(16, 2).Value = "3.8% (n = 18)"
$ws1.Cells.Item(16, 3).Value = "2.8% (n = 12)"
$ws1.Cells.Item(16, 4).Value = "ns (p = 0.67)"
$ws1.Cells.Item(16, 5).Value = "V = 0.026"

$ws1.Cells.Item(17, 1).Value = "Freq. bact. Infections"
$ws1.Cells.Item(17, 2).Value = "4.8% (n = 23)"
$ws1.Cells.Item(17, 3).Value = "1.2% (n = 5)"
$ws1.Cells.Item(17, 4).Value = "p = 0.016"
$ws1.Cells.Item(17, 5).Value = "V = 0.1"

$ws1.Cells.Item(18, 1).Value = "Pre-CoV sleep disorders"
$ws1.Cells.Item(18, 2).Value = "3.5% (n = 17)"
$ws1.Cells.Item(18, 3).Value = "4.7% (n = 20)"
$ws1.Cells.Item(18, 4).Value = "ns (p = 0.62)"
$ws1.Cells.Item(18, 5).Value = "V = 0.029"

$ws1.Cells.Item(19, 1).Value = "Daily medication"
$ws1.Cells.Item(19, 2).Value = "absent: 62% (n = 295)`n1 - 4 drugs: 37% (n = 175)`n5 drugs and more: 1.9% (n = 9)"
$ws1.Cells.Item(19, 3).Value = "absent: 74% (n = 317)`n1 - 4 drugs: 25% (n = 106)`n5 drugs and more: 0.94% (n = 4)"
$ws1.Cells.Item(19, 4).Value = "p = 0.0024"
$ws1.Cells.Item(19, 5).Value = "V = 0.14"

$ws1.Cells.Item(20, 1).Value = "Observation time"
$ws1.Cells.Item(20, 2).Value = "180 [IQR: 130 - 220]`nrange: 90 - 400"
$ws1.Cells.Item(20, 3).Value = "140 [IQR: 120 - 270]`nrange: 90 - 390"
$ws1.Cells.Item(20, 4).Value = "p = 0.0036"
$ws1.Cells.Item(20, 5).Value = "r = 0.12"

$ws1.Cells.Item(21, 1).Value = "Comorbidity"
$ws1.Cells.Item(21, 2).Value = "49% (n = 237)"
$ws1.Cells.Item(21, 3).Value = "43% (n = 185)"
$ws1.Cells.Item(21, 4).Value = "ns (p = 0.22)"
$ws1.Cells.Item(21, 5).Value = "V = 0.062"

# ---- Table 2 (Sheet2): reorder rows, drop "complete: n = ..." lines ----
$ws2 = $wb.Worksheets.Item(2)

$ws2.Cells.Item(2, 1).Value = "Sex"
$ws2.Cells.Item(2, 2).Value = "female: 41% (n = 44)`nmale: 59% (n = 64)"
$ws2.Cells.Item(2, 3).Value = "female: 67% (n = 18)`nmale: 33% (n = 9)"
$ws2.Cells.Item(2, 4).Value = "female: 35% (n = 19)`nmale: 65% (n = 36)"
$ws2.Cells.Item(2, 5).Value = "female: 27% (n = 7)`nmale: 73% (n = 19)"
$ws2.Cells.Item(2, 6).Value = "p < 0.001"
$ws2.Cells.Item(2, 7).Value = "V = 0.31"

$ws2.Cells.Item(3, 1).Value = "Age, years"
$ws2.Cells.Item(3, 2).Value = "56 [IQR: 49 - 68]`nrange: 19 - 87"
$ws2.Cells.Item(3, 3).Value = "47 [IQR: 38 - 55]`nrange: 19 - 70"
$ws2.Cells.Item(3, 4).Value = "62 [IQR: 53 - 72]`nrange: 27 - 87"
$ws2.Cells.Item(3, 5).Value = "56 [IQR: 52 - 64]`nrange: 44 - 79"
$ws2.Cells.Item(3, 6).Value = "p < 0.001"
$ws2.Cells.Item(3, 7).Value = "η² = 0.21"

$ws2.Cells.Item(4, 1).Value = "BMI at CoV onset"
$ws2.Cells.Item(4, 2).Value = "normal: 39% (n = 42)`noverweight: 43% (n = 46)`nobesity: 19% (n = 20)"
$ws2.Cells.Item(4, 3).Value = "normal: 56% (n = 15)`noverweight: 33% (n = 9)`nobesity: 11% (n = 3)"
$ws2.Cells.Item(4, 4).Value = "normal: 29% (n = 16)`noverweight: 51% (n = 28)`nobesity: 20% (n = 11)"
$ws2.Cells.Item(4, 5).Value = "normal: 42% (n = 11)`noverweight: 35% (n = 9)`nobesity: 23% (n = 6)"
$ws2.Cells.Item(4, 6).Value = "p < 0.001"
$ws2.Cells.Item(4, 7).Value = "V = 0.17"

$ws2.Cells.Item(5, 1).Value = "Comorbidity present"
$ws2.Cells.Item(5, 2).Value = "75% (n = 81)"
$ws2.Cells.Item(5, 3).Value = "41% (n = 11)"
$ws2.Cells.Item(5, 4).Value = "85% (n = 47)"
$ws2.Cells.Item(5, 5).Value = "88% (n = 23)"
$ws2.Cells.Item(5, 6).Value = "p < 0.001"
$ws2.Cells.Item(5, 7).Value = "V = 0.46"

$ws2.Cells.Item(6, 1).Value = "Cardiovascular disease"
$ws2.Cells.Item(6, 2).Value = "40% (n = 43)"
$ws2.Cells.Item(6, 3).Value = "7.4% (n = 2)"
$ws2.Cells.Item(6, 4).Value = "47% (n = 26)"
$ws2.Cells.Item(6, 5).Value = "58% (n = 15)"
$ws2.Cells.Item(6, 6).Value = "p < 0.001"
$ws2.Cells.Item(6, 7).Value = "V = 0.39"

$ws2.Cells.Item(7, 1).Value = "Hypertension"
$ws2.Cells.Item(7, 2).Value = "27% (n = 29)"
$ws2.Cells.Item(7, 3).Value = "7.4% (n = 2)"
$ws2.Cells.Item(7, 4).Value = "27% (n = 15)"
$ws2.Cells.Item(7, 5).Value = "46% (n = 12)"
$ws2.Cells.Item(7, 6).Value = "p < 0.001"
$ws2.Cells.Item(7, 7).Value = "V = 0.31"

$ws2.Cells.Item(8, 1).Value = "Pulmonary disease"
$ws2.Cells.Item(8, 2).Value = "19% (n = 20)"
$ws2.Cells.Item(8, 3).Value = "11% (n = 3)"
$ws2.Cells.Item(8, 4).Value = "22% (n = 12)"
$ws2.Cells.Item(8, 5).Value = "19% (n = 5)"
$ws2.Cells.Item(8, 6).Value = "p = 0.031"
$ws2.Cells.Item(8, 7).Value = "V = 0.11"

$ws2.Cells.Item(9, 1).Value = "Metabolic disease"
$ws2.Cells.Item(9, 2).Value = "42% (n = 45)"
$ws2.Cells.Item(9, 3).Value = "19% (n = 5)"
$ws2.Cells.Item(9, 4).Value = "49% (n = 27)"
$ws2.Cells.Item(9, 5).Value = "50% (n = 13)"
$ws2.Cells.Item(9, 6).Value = "p < 0.001"
$ws2.Cells.Item(9, 7).Value = "V = 0.27"

$ws2.Cells.Item(10, 1).Value = "Diabetes"
$ws2.Cells.Item(10, 2).Value = "15% (n = 16)"
$ws2.Cells.Item(10, 3).Value = "3.7% (n = 1)"
$ws2.Cells.Item(10, 4).Value = "15% (n = 8)"
$ws2.Cells.Item(10, 5).Value = "27% (n = 7)"
$ws2.Cells.Item(10, 6).Value = "p < 0.001"
$ws2.Cells.Item(10, 7).Value = "V = 0.23"

$ws2.Cells.Item(11, 1).Value = "Gastrointestinal disease"
$ws2.Cells.Item(11, 2).Value = "13% (n = 14)"
$ws2.Cells.Item(11, 3).Value = "0% (n = 0)"
$ws2.Cells.Item(11, 4).Value = "20% (n = 11)"
$ws2.Cells.Item(11, 5).Value = "12% (n = 3)"
$ws2.Cells.Item(11, 6).Value = "p < 0.001"
$ws2.Cells.Item(11, 7).Value = "V = 0.24"

$ws2.Cells.Item(12, 1).Value = "Malignancy"
$ws2.Cells.Item(12, 2).Value = "9.3% (n = 10)"
$ws2.Cells.Item(12, 3).Value = "3.7% (n = 1)"
$ws2.Cells.Item(12, 4).Value = "15% (n = 8)"
$ws2.Cells.Item(12, 5).Value = "3.8% (n = 1)"
$ws2.Cells.Item(12, 6).Value = "p < 0.001"
$ws2.Cells.Item(12, 7).Value = "V = 0.19"

$ws2.Cells.Item(13, 1).Value = "Immune deficiency"
$ws2.Cells.Item(13, 2).Value = "5.6% (n = 6)"
$ws2.Cells.Item(13, 3).Value = "0% (n = 0)"
$ws2.Cells.Item(13, 4).Value = "3.6% (n = 2)"
$ws2.Cells.Item(13, 5).Value = "15% (n = 4)"
$ws2.Cells.Item(13, 6).Value = "p < 0.001"
$ws2.Cells.Item(13, 7).Value = "V = 0.25"

